# Vers 1.1 BMC Gastroenterology - publication version
#
# Adds a new "NOTES" worksheet (placed after the existing data sheet) that
# documents a few of the abbreviations used on the main sheet, and makes the
# new sheet the active/selected tab.

$wb = $excel.ActiveWorkbook

# The existing (only) worksheet becomes the first sheet; add the new NOTES
# sheet directly after it so it keeps the sheet order: data sheet, NOTES.
$dataSheet = $wb.Worksheets.Item(1)
$notes = $wb.Worksheets.Add($null, $dataSheet)
$notes.Name = "NOTES"

$notes.Range("A1").Value = "ACM"
$notes.Range("B1").Value = "all cause mortality"
$notes.Range("A2").Value = "wm/bm/wf"
$notes.Range("B2").Value = "white or black; male or female"
$notes.Range("A3").Value = "mortality_xx_EA"
$notes.Range("B3").Value = "mortality from esophageal adenocarcinoma"

$notes.Columns.Item(1).AutoFit() | Out-Null
$notes.Columns.Item(2).AutoFit() | Out-Null

$notes.Range("A4").Select() | Out-Null
